$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 5 (pushes the old row 5 and everything below it down by
# one; formulas / shared-formula ranges auto-adjust).
$ws.Rows(5).Insert()

# Drop the old "EVP<2>" label that used to live at L15 (now shifted to L16).
$ws.Range("L16").ClearContents()

# Row 3: "neighbours =" / 8  ->  "num_neighbours_in_nn_table" / 12
$ws.Range("C3").Value2 = "num_neighbours_in_nn_table"
$ws.Range("D3").Value2 = 12

# New row 5: "build_reverse_list_size" / 32
$ws.Range("C5").Value2 = "build_reverse_list_size"
$ws.Range("D5").Value2 = 32

# Row 4 (used to be "rev_neighbours" / 20) -> "num_neighbours_in_reverse_table" / 16
$ws.Range("C4").Value2 = "num_neighbours_in_reverse_table"
$ws.Range("D4").Value2 = 16

# New little division/multiplication block at the bottom of the sheet.
$ws.Range("F21").Value2 = 384
$ws.Range("G21").Value2 = "/"
$ws.Range("H21").Value2 = 64
$ws.Range("I21").Formula = "=F21/H21"

$ws.Range("G22").Value2 = "*"
$ws.Range("H22").Formula = "=H21*I21"

# Rename "EVP<2>" -> "EVP" for the first data row (now row 10).
$ws.Range("B10").Value2 = "EVP"

# Widen column C to fit the longer labels now in use.
$ws.Columns(3).ColumnWidth = 40.6

# Move the active selection, matching the saved view state.
$ws.Range("F4").Select() | Out-Null
